# "Generate Report for Handback": the bef00e20-... file has now been
# handed back and is in sync with en-US, so refresh its status/date/error
# fields across the Overview, zh-cn and de-de report sheets (row 3 in
# each sheet is the bef00e20-... file).

$wb = $excel.ActiveWorkbook

# --- Overview sheet: Status columns (zh-cn / de-de) for the handed-back file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn detail sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K3").Value = "2016-08-25 20:48:17"
$wsZhCn.Range("P3").Value = ""
$wsZhCn.Columns.Item(16).AutoFit()

# --- de-de detail sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K3").Value = "2016-08-25 20:48:24"
$wsDeDe.Range("P3").Value = ""
$wsDeDe.Columns.Item(16).AutoFit()
